# Generate Report for Handback
# Refresh the "Latest HO Xliff Generate Date" / handback timestamp cells on
# each sheet to reflect the new report-generation run.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date for 4a0a440b-...md
$overview.Range("G4").Value = "2016-08-15 09:00:55"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
# for 4a0a440b-....zh-cn.xlf
$zhcn.Range("H4").Value = "2016-08-15 09:00:50"
$zhcn.Range("K4").Value = "2016-08-15 09:01:14"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime
# for 4a0a440b-....de-de.xlf
$dede.Range("H4").Value = "2016-08-15 09:00:55"
$dede.Range("K4").Value = "2016-08-15 09:01:21"
